# Update cryptos list data (prices and volume %), and a few rows where
# coin rankings changed position / a coin was replaced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.486.72"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "3.358.34"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D5").Value = "258.25"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "664.69"
$ws.Range("E6").Value = "  +5.23%  "
$ws.Range("D7").Value = "1.53"
$ws.Range("E7").Value = "  +9.02%  "
$ws.Range("D8").Value = "0.467"
$ws.Range("E8").Value = "  +19.99%  "
$ws.Range("D9").Value = "1.07"
$ws.Range("E9").Value = "  +24.39%  "
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").Value = "3.352.43"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "0.213"
$ws.Range("E12").Value = "  +7.33%  "
$ws.Range("D13").Value = "42.21"
$ws.Range("E13").Value = "  +15.35%  "
$ws.Range("E14").Value = "  +9.57%  "
$ws.Range("D15").Value = "98.704.31"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "3.992.41"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "5.64"
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").Value = "3.362.41"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  +25.55%  "
$ws.Range("D20").Value = "16.85"
$ws.Range("E20").Value = "  +10.32%  "
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").Value = "532.76"
$ws.Range("E22").Value = "  +8.77%  "
$ws.Range("D23").Value = "10.50"
$ws.Range("E23").Value = "  +10.58%  "
$ws.Range("D24").Value = "0.0000218"
$ws.Range("E24").Value = "  +3.81%  "
$ws.Range("D25").Value = "0.434"
$ws.Range("E25").Value = "  +52.40%  "
$ws.Range("D26").Value = "102.12"
$ws.Range("E26").Value = "  +14.59%  "
$ws.Range("D27").Value = "6.23"
$ws.Range("E27").Value = "  +10.32%  "
$ws.Range("D28").Value = "12.58"
$ws.Range("E28").Value = "  +5.48%  "
$ws.Range("D29").Value = "3.541.51"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "0.148"
$ws.Range("E30").Value = "  +8.54%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "10.99"
$ws.Range("E32").Value = "  +13.19%  "
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D35").Value = "29.40"
$ws.Range("E35").Value = "  +4.52%  "
$ws.Range("D36").Value = "0.544"
$ws.Range("E36").Value = "  +17.85%  "
$ws.Range("D37").Value = "7.83"
$ws.Range("E37").Value = "  +7.24%  "
$ws.Range("D38").Value = "2.12"
$ws.Range("E38").Value = "  +8.16%  "
$ws.Range("E39").Value = "  +4.75%  "
$ws.Range("D40").Value = "525.89"
$ws.Range("E40").Value = "  +5.24%  "
$ws.Range("D41").Value = "1.34"
$ws.Range("E41").Value = "  +5.96%  "
$ws.Range("D42").Value = "24.69"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").Value = "0.0439"
$ws.Range("E43").Value = "  +33.75%  "
$ws.Range("D44").Value = "3.80"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "3.41"
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "0.829"
$ws.Range("E46").Value = "  +4.82%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "2.08"
$ws.Range("E48").Value = "  +7.22%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "7.95"
$ws.Range("E49").Value = "  +19.05%  "
$ws.Range("D50").Value = "5.16"
$ws.Range("E50").Value = "  +11.53%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "50.76"
$ws.Range("E51").Value = "  +10.66%  "
